$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Move-RowCells($srcRow, $dstRow, $cols) {
    foreach ($col in $cols) {
        $ws.Cells.Item($srcRow, $col).Copy($ws.Cells.Item($dstRow, $col))
    }
}

# ------------------------------------------------------------------
# 1) Push the existing header/body block (old rows 6-10) down by 6
#    rows so it lands on rows 12-16. Work bottom-up so we never
#    clobber a row before it has been copied. Row 11 is intentionally
#    never touched/written, so it stays absent from the sheet, just
#    like the old spacer row 5 vanishes instead of shifting into it.
# ------------------------------------------------------------------
Move-RowCells 10 16 @(1,6,7,8,9,10,11)
Move-RowCells 9  15 @(1)
Move-RowCells 8  14 @(1,2,3,4,5,6,7,8,9,10,11)
Move-RowCells 7  13 @(1,4,5)
Move-RowCells 6  12 @(1,2,3,4,5,6,7,8,9,10,11)

# Row 16 (old totals row) keeps the sheet's default height (no
# explicit override existed before the move either).
$ws.Rows(15).RowHeight = 12.1
$ws.Rows(14).RowHeight = 12.1
$ws.Rows(13).RowHeight = 12.2
$ws.Rows(12).RowHeight = 12.2

# ------------------------------------------------------------------
# 2) Build the new tracking rows 5-10 using row 4's look (maroon text
#    on the lavender band), each merged A:K like the banner rows
#    above them. ClearFormats first so no stale row-level style/
#    customFormat flag survives from the header content that used to
#    live in rows 6-10.
# ------------------------------------------------------------------
$ws.Rows("5:10").ClearFormats()

$trackingRows = @(
    @{ Row = 5;  Text = "{{#if track_name}}" },
    @{ Row = 6;  Text = "Tracking: {{track_name}}" },
    @{ Row = 7;  Text = "{{/if}}" },
    @{ Row = 8;  Text = "{{#if track2_name}}" },
    @{ Row = 9;  Text = "Tracking-2: {{track2_name}}" },
    @{ Row = 10; Text = "{{/if}}" }
)

foreach ($item in $trackingRows) {
    $r = $item.Row
    for ($col = 1; $col -le 11; $col++) {
        $ws.Cells.Item(4, $col).Copy($ws.Cells.Item($r, $col))
    }
    $ws.Rows($r).RowHeight = 12.2
    $ws.Cells.Item($r, 1).Value = $item.Text
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 11)).Merge()
}

# ------------------------------------------------------------------
# 3) Fix up the view: selection back to A8, scrolled to the top-left.
# ------------------------------------------------------------------
$ws.Range("A8").Select()

Write-Output ("Dimension: " + $ws.UsedRange.Address())
